# Apply odds updates to row 5 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 2.3

$ws.Range("AA5").Value = 15
$ws.Range("AD5").Value = 7
$ws.Range("AF5").Value = 81

$ws.Range("AN5").Value = 3.4
$ws.Range("AQ5").Value = 29
